$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: for numeric-looking text values (e.g. "92.38"), Excel
# would normally auto-convert the assigned string into a Number cell. The source
# workbook stores these Price/Volume columns as literal text, so we force the
# "Text" number format just for the assignment, then clear the temporary format
# so the cell keeps its original (default) style while remaining text-typed.

$ws.Range('D2').Value = '30.188.62'
$ws.Range('E2').Value = '  -0.21%  '

$ws.Range('D3').Value = '1.899.24'
$ws.Range('E3').Value = '  -1.12%  '

$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range('E4').Value = '  +0.12%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '321.11'
$c.ClearFormats()
$ws.Range('E5').Value = '  -2.82%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range('E6').Value = '  +0.16%  '

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.5055'
$c.ClearFormats()
$ws.Range('E7').Value = '  -3.29%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.4033'
$c.ClearFormats()
$ws.Range('E8').Value = '  -1.28%  '

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.08277'
$c.ClearFormats()
$ws.Range('E9').Value = '  -3.03%  '

$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '1.110'
$c.ClearFormats()
$ws.Range('E10').Value = '  -1.60%  '

$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '42.28'
$c.ClearFormats()
$ws.Range('E11').Value = '  -1.46%  '

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '24.15'
$c.ClearFormats()
$ws.Range('E12').Value = '  +5.82%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '6.404'
$c.ClearFormats()
$ws.Range('E13').Value = '  -0.73%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.883.63'
$ws.Range('E14').Value = '  -1.57%  '

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '7.315'
$c.ClearFormats()
$ws.Range('E15').Value = '  -1.36%  '

$ws.Range('E16').Value = '  +0.42%  '

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '92.38'
$c.ClearFormats()
$ws.Range('E17').Value = '  -3.01%  '

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.00001098'
$c.ClearFormats()
$ws.Range('E18').Value = '  -1.45%  '

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.06472'
$c.ClearFormats()
$ws.Range('E19').Value = '  -3.36%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '18.45'
$c.ClearFormats()
$ws.Range('E20').Value = '  +0.17%  '

$ws.Range('E21').Value = '  +0.18%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '5.932'
$c.ClearFormats()
$ws.Range('E22').Value = '  -1.43%  '

$ws.Range('D23').Value = '30.198.11'
$ws.Range('E23').Value = '  -0.16%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '11.27'
$c.ClearFormats()
$ws.Range('E24').Value = '  -0.66%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.193'
$c.ClearFormats()
$ws.Range('E25').Value = '  -1.31%  '

$ws.Range('D26').Value = '2.103.20'
$ws.Range('E26').Value = '  -1.57%  '

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '21.58'
$c.ClearFormats()
$ws.Range('E27').Value = '  +2.04%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '160.88'
$c.ClearFormats()
$ws.Range('E28').Value = '  -0.09%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '2.298'
$c.ClearFormats()
$ws.Range('E29').Value = '  -5.38%  '

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '129.12'
$c.ClearFormats()
$ws.Range('E30').Value = '  +0.02%  '

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.123'
$c.ClearFormats()
$ws.Range('E31').Value = '  +3.38%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.1042'
$c.ClearFormats()
$ws.Range('E32').Value = '  -2.14%  '

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '5.979'
$c.ClearFormats()
$ws.Range('E33').Value = '  -1.12%  '

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '3.690'
$c.ClearFormats()
$ws.Range('E34').Value = '  +1.98%  '

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.02445'
$c.ClearFormats()
$ws.Range('E35').Value = '  -1.91%  '

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '5.322'
$c.ClearFormats()
$ws.Range('E36').Value = '  +2.61%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.06450'
$c.ClearFormats()
$ws.Range('E37').Value = '  -1.95%  '

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.2155'
$c.ClearFormats()
$ws.Range('E38').Value = '  -2.55%  '

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.186'
$c.ClearFormats()
$ws.Range('E39').Value = '  -3.73%  '

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '8.608'
$c.ClearFormats()
$ws.Range('E40').Value = '  -2.96%  '

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.6369'
$c.ClearFormats()
$ws.Range('E41').Value = '  -2.77%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '11.36'
$c.ClearFormats()
$ws.Range('E42').Value = '  -3.16%  '

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.214'
$c.ClearFormats()
$ws.Range('E43').Value = '  -2.36%  '

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range('E44').Value = '  +0.32%  '

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '13.23'
$c.ClearFormats()
$ws.Range('E45').Value = '  +0.04%  '

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.5980'
$c.ClearFormats()
$ws.Range('E46').Value = '  -2.97%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '2.157'
$c.ClearFormats()
$ws.Range('E47').Value = '  +3.21%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '3.638'
$c.ClearFormats()
$ws.Range('E48').Value = '  -3.10%  '

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '123.53'
$c.ClearFormats()
$ws.Range('E49').Value = '  -0.56%  '

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.215'
$c.ClearFormats()
$ws.Range('E50').Value = '  -2.69%  '

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '78.60'
$c.ClearFormats()
$ws.Range('E51').Value = '  -1.40%  '
